$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 902.10345
$ws.Range("J17").Value = 898.9643
$ws.Range("L17").Value = 2696.8929
$ws.Range("N17").Value = -3032.8929
$ws.Range("H38").Value = 2370.8667
$ws.Range("J38").Value = 3371.2856
$ws.Range("L38").Value = 10113.8568
$ws.Range("N38").Value = -10857.8568
$ws.Range("H42").Value = 360.16666
$ws.Range("J42").Value = 759.4
$ws.Range("L42").Value = 2278.2
$ws.Range("N42").Value = -2738.2
$ws.Range("H46").Value = 168763.17
$ws.Range("I46").Value = 793
$ws.Range("J46").Value = 336733.34
$ws.Range("K46").Value = 2379
$ws.Range("L46").Value = 1010200.02
$ws.Range("M46").Value = -2260
$ws.Range("N46").Value = -1010438.02
$ws.Range("H58").Value = 2251.5386
$ws.Range("J58").Value = 3702.4285
$ws.Range("L58").Value = 11107.2855
$ws.Range("N58").Value = -11407.2855
$ws.Range("H60").Value = 168763.17
$ws.Range("I60").Value = 793
$ws.Range("J60").Value = 336733.34
$ws.Range("K60").Value = 2379
$ws.Range("L60").Value = 1010200.02
$ws.Range("M60").Value = -1895
$ws.Range("N60").Value = -1011168.02
$ws.Range("H61").Value = 807.6667
$ws.Range("I61").Value = 807.6667
$ws.Range("K61").Value = 2423.0001
$ws.Range("M61").Value = -2251.0001
$ws.Range("H99").Value = 1642332
$ws.Range("I99").Value = 2343248.8
$ws.Range("J99").Value = 6859.3335
$ws.Range("K99").Value = 7029746.399999999
$ws.Range("L99").Value = 20578.0005
$ws.Range("M99").Value = -7028248.399999999
$ws.Range("N99").Value = -23574.0005
$ws.Range("H101").Value = 12989041
$ws.Range("I101").Value = 20409992
$ws.Range("K101").Value = 61229976
$ws.Range("M101").Value = -61228354
$ws.Range("H111").Value = 3510.6667
$ws.Range("I111").Value = 3510.6667
$ws.Range("K111").Value = 10532.0001
$ws.Range("M111").Value = -7465.000100000001
$ws.Range("H112").Value = 54467.473
$ws.Range("J112").Value = 68639.07000000001
$ws.Range("L112").Value = 205917.21
$ws.Range("N112").Value = -208133.21
$ws.Range("H115").Value = 735.5
$ws.Range("I115").Value = 412.14285
$ws.Range("K115").Value = 1236.42855
$ws.Range("M115").Value = 330.5714499999999
$ws.Range("H127").Value = 5170.778
$ws.Range("I127").Value = 1770.1666
$ws.Range("J127").Value = 11972
$ws.Range("K127").Value = 5310.4998
$ws.Range("L127").Value = 35916
$ws.Range("M127").Value = -350.4997999999996
$ws.Range("N127").Value = -45836
$ws.Range("H129").Value = 1829.6666
$ws.Range("I129").Value = 1126.3334
$ws.Range("J129").Value = 2533
$ws.Range("K129").Value = 3379.0002
$ws.Range("L129").Value = 7599
$ws.Range("M129").Value = 1620.9998
$ws.Range("N129").Value = -17599
$ws.Range("H131").Value = 5818
$ws.Range("I131").Value = 2000
$ws.Range("K131").Value = 6000
$ws.Range("M131").Value = -960
$ws.Range("H135").Value = 3362.4666
$ws.Range("I135").Value = 3362.4666
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 30262.1994
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -27727.1994
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 3423.6038
$ws.Range("I138").Value = 1490.4375
$ws.Range("J138").Value = 4259.5674
$ws.Range("K138").Value = 4471.3125
$ws.Range("L138").Value = 12778.7022
$ws.Range("M138").Value = 668.6875
$ws.Range("N138").Value = -23058.7022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -598
$ws.Range("H32").Value = 1941.0322
$ws.Range("I32").Value = 1956.678
$ws.Range("K32").Value = 1956.678
$ws.Range("M32").Value = -1669.678
$ws.Range("H35").Value = 11999
$ws.Range("I35").Value = 11999
$ws.Range("K35").Value = 11999
$ws.Range("M35").Value = -11593
$ws.Range("H45").Value = 7249.846
$ws.Range("J45").Value = 6581
$ws.Range("L45").Value = 6581
$ws.Range("N45").Value = -7335
$ws.Range("H74").Value = 3114.4167
$ws.Range("I74").Value = 1736.6316
$ws.Range("K74").Value = 1736.6316
$ws.Range("M74").Value = -862.6315999999999
$ws.Range("H77").Value = 3114.4167
$ws.Range("I77").Value = 1736.6316
$ws.Range("K77").Value = 8683.157999999999
$ws.Range("M77").Value = -4315.157999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 5448.5
$ws.Range("I25").Value = 5399.5
$ws.Range("J25").Value = 5497.5
$ws.Range("K25").Value = 5399.5
$ws.Range("L25").Value = 5497.5
$ws.Range("M25").Value = -5164.5
$ws.Range("N25").Value = -5967.5
$ws.Range("H37").Value = 4875
$ws.Range("I37").Value = 2812.5
$ws.Range("K37").Value = 2812.5
$ws.Range("M37").Value = -2675.5
$ws.Range("H56").Value = 24666.666
$ws.Range("I56").Value = 45000
$ws.Range("K56").Value = 45000
$ws.Range("M56").Value = -44261
$ws.Range("H86").Value = 7644.3477
$ws.Range("I86").Value = 6352.8125
$ws.Range("K86").Value = 6352.8125
$ws.Range("M86").Value = -5229.8125
$ws.Range("H89").Value = 7644.3477
$ws.Range("I89").Value = 6352.8125
$ws.Range("K89").Value = 31764.0625
$ws.Range("M89").Value = -26148.0625
$ws.Range("H107").Value = 3595.76
$ws.Range("I107").Value = 3305.5293
$ws.Range("K107").Value = 3305.5293
$ws.Range("M107").Value = -1385.5293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2780.5386
$ws.Range("I31").Value = 1711.1
$ws.Range("J31").Value = 3906.2632
$ws.Range("K31").Value = 1711.1
$ws.Range("L31").Value = 3906.2632
$ws.Range("M31").Value = -1416.1
$ws.Range("N31").Value = -4496.263199999999
$ws.Range("H34").Value = 2780.5386
$ws.Range("I34").Value = 1711.1
$ws.Range("J34").Value = 3906.2632
$ws.Range("K34").Value = 1711.1
$ws.Range("L34").Value = 3906.2632
$ws.Range("M34").Value = -1509.1
$ws.Range("N34").Value = -4310.263199999999
$ws.Range("H58").Value = 1718.1428
$ws.Range("I58").Value = 940.931
$ws.Range("K58").Value = 940.931
$ws.Range("M58").Value = -737.931
$ws.Range("H134").Value = 3586.9736
$ws.Range("I134").Value = 1718.7693
$ws.Range("K134").Value = 5156.3079
$ws.Range("M134").Value = -2621.3079
$ws.Range("H136").Value = 1718.1428
$ws.Range("I136").Value = 940.931
$ws.Range("K136").Value = 2822.793
$ws.Range("M136").Value = -272.7930000000001
$ws.Range("H139").Value = 55748.332
$ws.Range("J139").Value = 55748.332
$ws.Range("L139").Value = 55748.332
$ws.Range("N139").Value = -66028.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 324099.84
$ws.Range("I5").Value = 1040.6666
$ws.Range("K5").Value = 3121.9998
$ws.Range("M5").Value = -3009.9998
$ws.Range("H68").Value = 35717948
$ws.Range("I68").Value = 2749.6667
$ws.Range("K68").Value = 8249.000100000001
$ws.Range("M68").Value = -7438.000100000001
$ws.Range("H71").Value = 35717948
$ws.Range("I71").Value = 2749.6667
$ws.Range("K71").Value = 24747.0003
$ws.Range("M71").Value = -20691.0003
$ws.Range("H135").Value = 324099.84
$ws.Range("I135").Value = 1040.6666
$ws.Range("K135").Value = 9365.999400000001
$ws.Range("M135").Value = -6830.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3385.6843
$ws.Range("J132").Value = 2471.1667
$ws.Range("L132").Value = 7413.500100000001
$ws.Range("N132").Value = -12473.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25671.182
$ws.Range("I40").Value = 27543.111
$ws.Range("K40").Value = 27543.111
$ws.Range("M40").Value = -27407.111
$ws.Range("H46").Value = 2830.087
$ws.Range("I46").Value = 1845.2727
$ws.Range("J46").Value = 3732.8333
$ws.Range("K46").Value = 1845.2727
$ws.Range("L46").Value = 3732.8333
$ws.Range("M46").Value = -1657.2727
$ws.Range("N46").Value = -4108.8333
$ws.Range("H127").Value = 250180000
$ws.Range("J127").Value = 240000
$ws.Range("L127").Value = 240000
$ws.Range("N127").Value = -249920
$ws.Range("H136").Value = 5344.5806
$ws.Range("I136").Value = 1825.3077
$ws.Range("K136").Value = 5475.9231
$ws.Range("M136").Value = -2925.9231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 73667.336
$ws.Range("I3").Value = 73667.336
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 73667.336
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -73553.336
$ws.Range("N3").ClearContents()
$ws.Range("H51").Value = 59994.668
$ws.Range("I51").Value = 59994
$ws.Range("J51").Value = 59995
$ws.Range("K51").Value = 59994
$ws.Range("L51").Value = 59995
$ws.Range("M51").Value = -59484
$ws.Range("N51").Value = -61015
$ws.Range("H61").Value = 30944.75
$ws.Range("I61").Value = 31266.334
$ws.Range("K61").Value = 31266.334
$ws.Range("M61").Value = -30974.334
$ws.Range("H81").Value = 13853.923
$ws.Range("I81").Value = 14341.75
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 28683.5
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = -27622.5
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 13853.923
$ws.Range("I84").Value = 14341.75
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 143417.5
$ws.Range("L84").Value = 80000
$ws.Range("M84").Value = -138113.5
$ws.Range("N84").Value = -90608
$ws.Range("H107").Value = 156249.75
$ws.Range("I107").Value = 12500
$ws.Range("K107").Value = 37500
$ws.Range("M107").Value = -35580
$ws.Range("H136").Value = 1951.2858
$ws.Range("I136").Value = 1524.4615
$ws.Range("K136").Value = 4573.3845
$ws.Range("M136").Value = -2023.3845
